$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-13 Sunday", "2025-04-14 Monday"),
    @("949÷4=", "396÷4="),
    @("121÷5=", "274÷2="),
    @("564÷7=", "360÷4="),
    @("102÷2=", "921÷9="),
    @("902÷9=", "315÷7="),
    @("747÷2=", "734÷6="),
    @("570÷7=", "582÷8="),
    @("228÷2=", "515÷2="),
    @("901÷5=", "914÷3="),
    @("942÷2=", "270÷5="),
    @("311÷9=", "557÷5="),
    @("127÷4=", "755÷9="),
    @("546÷3=", "468÷2="),
    @("169÷9=", "220÷6="),
    @("487÷5=", "573÷7="),
    @("516÷8=", "993÷2="),
    @("425÷9=", "930÷6="),
    @("506÷2=", "680÷6="),
    @("299÷8=", "930÷4="),
    @("758÷4=", "478÷2="),
    @("951÷8=", "790÷5="),
    @("406÷9=", "491÷8="),
    @("343÷3=", "723÷4="),
    @("223÷4=", "380÷6="),
    @("728÷6=", "759÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
